$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(0, 3, 0, 3, 2, 8)
    3 = @(6, 5, 6, 3, 9, 6)
    4 = @(0, 9, 6, 7, 7, 0)
    5 = @(7, 3, 7, 1, 2, 1)
    6 = @(3, 5, 1, 3, 7, 6)
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = [char](68 + $i)  # 'D' is 68
        $ws.Range("$colLetter$row").Value = $cols[$i]
    }
}
